# Updates cryptocurrency price (D) and 1h volume-change (E) columns
# to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.966.70"
$ws.Range("E2").Value = "  +1.21%  "

$ws.Range("D3").Value = "1.719.83"
$ws.Range("E3").Value = "  +1.42%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3976"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.69%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4127"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.64%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.002"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.28%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.63%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08985"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.713"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.70%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.214"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.13%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001374"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.48%  "

$ws.Range("D17").Value = "1.702.36"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "100.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07159"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.57%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.522"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.48%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.32%  "

$ws.Range("D24").Value = "24.962.43"
$ws.Range("E24").Value = "  +1.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.153"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.342"
$ws.Range("D26").Style = "Normal"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.305"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +25.18%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.43"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "139.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.238"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.892"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.20%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09113"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.54%  "

$ws.Range("D34").Value = "1.888.46"
$ws.Range("E34").Value = "  +0.02%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.100"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.61%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03024"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2827"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.20"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.966"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "14.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.09355"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8147"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.31%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.496"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.14%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7412"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.655"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.277"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.358"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.56%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "94.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.00%  "
